# Actualización automática del tracker
# - Completa el resultado de la fila 12 (G12/H12)
# - Añade dos nuevas filas de partidos (21 y 22) al final de la tabla

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fila 12: completar resultado pendiente ---
$ws.Range("G12").Value = "Fallo"
$ws.Range("H12").Value = -1

# --- Nuevas filas de partidos ---
$newRows = @(
    @{ Row = 21; EventId = 14266314; Fecha = "2025-08-02"; JugadorA = "Elena Rybakina"; JugadorB = "Dayana Yastremska"; Pronostico = "Gana Elena Rybakina"; Cuota = 2.3;  Resultado = "Acierto"; Profit = 1.3 },
    @{ Row = 22; EventId = 14266308; Fecha = "2025-08-02"; JugadorA = "Coco Gauff";     JugadorB = "Victoria Mboko";      Pronostico = "Gana Coco Gauff";      Cuota = 1.4;  Resultado = "Fallo";   Profit = -1  }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.EventId

    # La fecha debe guardarse como texto literal (no como fecha de Excel),
    # igual que el resto de la columna "fecha".
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $r.Fecha
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").Value = $r.JugadorA
    $ws.Range("D$row").Value = $r.JugadorB
    $ws.Range("E$row").Value = $r.Pronostico
    $ws.Range("F$row").Value = $r.Cuota
    $ws.Range("G$row").Value = $r.Resultado
    $ws.Range("H$row").Value = $r.Profit
}
